$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.677.62"
$ws.Range("E2").Value = "  -6.56%  "
$ws.Range("D3").Value = "'3.265.92"
$ws.Range("E3").Value = "  -9.10%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'176.98"
$ws.Range("E5").Value = "  -13.57%  "
$ws.Range("D6").Value = "'511.36"
$ws.Range("E6").Value = "  -10.44%  "
$ws.Range("D7").Value = "'0.585"
$ws.Range("E7").Value = "  -4.55%  "
$ws.Range("D8").Value = "'3.265.11"
$ws.Range("E8").Value = "  -9.02%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'0.608"
$ws.Range("E10").Value = "  -10.60%  "
$ws.Range("D11").Value = "'56.73"
$ws.Range("E11").Value = "  -6.68%  "
$ws.Range("E12").Value = "  -12.75%  "
$ws.Range("E13").Value = "  -10.72%  "
$ws.Range("E14").Value = "  -13.13%  "
$ws.Range("D15").Value = "'3.786.73"
$ws.Range("E15").Value = "  -8.92%  "
$ws.Range("E16").Value = "  -5.06%  "
$ws.Range("D17").Value = "'3.268.97"
$ws.Range("E17").Value = "  -8.85%  "
$ws.Range("D18").Value = "'63.498.24"
$ws.Range("E18").Value = "  -6.52%  "
$ws.Range("E19").Value = "  -11.00%  "
$ws.Range("D20").Value = "'10.68"
$ws.Range("E20").Value = "  -12.93%  "
$ws.Range("D21").Value = "'0.933"
$ws.Range("E21").Value = "  -12.13%  "
$ws.Range("D22").Value = "'365.28"
$ws.Range("E22").Value = "  -9.30%  "
$ws.Range("E23").Value = "  -6.64%  "
$ws.Range("E24").Value = "  -13.84%  "
$ws.Range("D25").Value = "'10.63"
$ws.Range("E25").Value = "  -14.90%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'3.77"
$ws.Range("E26").Value = "  -3.20%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'5.99"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("E28").Value = "  -10.17%  "
$ws.Range("D29").Value = "'11.13"
$ws.Range("E29").Value = "  -10.75%  "
$ws.Range("D30").Value = "'8.20"
$ws.Range("E30").Value = "  -11.28%  "
$ws.Range("D31").Value = "'28.03"
$ws.Range("E31").Value = "  -11.29%  "
$ws.Range("D32").Value = "'633.02"
$ws.Range("E32").Value = "  -5.63%  "
$ws.Range("D33").Value = "'6.58"
$ws.Range("E33").Value = "  -15.35%  "
$ws.Range("D34").Value = "'10.92"
$ws.Range("E34").Value = "  -9.82%  "
$ws.Range("D35").Value = "'58.47"
$ws.Range("E35").Value = "  -7.72%  "
$ws.Range("E36").Value = "  -10.30%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'35.15"
$ws.Range("E38").Value = "  -14.68%  "
$ws.Range("E39").Value = "  -10.21%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "'0.120"
$ws.Range("E41").Value = "  -9.49%  "
$ws.Range("D42").Value = "'2.825.98"
$ws.Range("E42").Value = "  -11.52%  "
$ws.Range("D43").Value = "'2.64"
$ws.Range("E43").Value = "  -19.27%  "
$ws.Range("D44").Value = "0.0₃0622"
$ws.Range("E44").Value = "  -17.48%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'25.90"
$ws.Range("E45").Value = "  +17.86%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'2.56"
$ws.Range("E46").Value = "  -8.70%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0379"
$ws.Range("E47").Value = "  -7.71%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.25"
$ws.Range("E48").Value = "  -16.44%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.123"
$ws.Range("E49").Value = "  -6.16%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'132.35"
$ws.Range("E50").Value = "  -4.82%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'2.63"
$ws.Range("E51").Value = "  -2.47%  "
